$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 45.31778466666666
$ws.Range("H2").Value = 135.953354
$ws.Range("I2").Value = 0.1102361023838286
$ws.Range("J2").Value = 0.1102361023838286
$ws.Range("M2").Value = 47.32925566666668
$ws.Range("N2").Value = 141.987767
$ws.Range("O2").Value = 0.3408416299313156
$ws.Range("P2").Value = 0.3408416299313156
$ws.Range("Q2").Value = 2144.857016735613
$ws.Range("R2").Value = 19303.71315062052
$ws.Range("S2").Value = 0.03757305281377953
$ws.Range("T2").Value = 0.03757305281377954
# Row 3
$ws.Range("G3").Value = 45.31778466666666
$ws.Range("H3").Value = 135.953354
$ws.Range("I3").Value = 0.1102361023838286
$ws.Range("J3").Value = 0.1102361023838286
$ws.Range("M3").Value = 43.717953
$ws.Range("O3").Value = 0.3148348341399153
$ws.Range("P3").Value = 0.3148348341399154
$ws.Range("Q3").Value = 1981.200780121454
$ws.Range("R3").Value = 17830.80702109309
$ws.Range("S3").Value = 0.03470616501024341
$ws.Range("T3").Value = 0.03470616501024342
# Row 4
$ws.Range("G4").Value = 45.31778466666666
$ws.Range("H4").Value = 135.953354
$ws.Range("I4").Value = 0.1102361023838286
$ws.Range("J4").Value = 0.1102361023838286
$ws.Range("M4").Value = 21.09134933333333
$ws.Range("N4").Value = 63.274048
$ws.Range("O4").Value = 0.1518893501062827
$ws.Range("P4").Value = 0.1518893501062827
$ws.Range("Q4").Value = 955.8132274174435
$ws.Range("R4").Value = 8602.319046756991
$ws.Range("S4").Value = 0.01674368994932937
$ws.Range("T4").Value = 0.01674368994932938
# Row 5
$ws.Range("G5").Value = 45.31778466666666
$ws.Range("H5").Value = 135.953354
$ws.Range("I5").Value = 0.1102361023838286
$ws.Range("J5").Value = 0.1102361023838286
$ws.Range("M5").Value = 26.72140366666666
$ws.Range("N5").Value = 80.16421099999999
$ws.Range("O5").Value = 0.1924341858224863
$ws.Range("P5").Value = 0.1924341858224864
$ws.Range("Q5").Value = 1210.954817357077
$ws.Range("R5").Value = 10898.59335621369
$ws.Range("S5").Value = 0.02121319461047631
$ws.Range("T5").Value = 0.02121319461047631
# Row 6
$ws.Range("I6").Value = 0.2429203181515272
$ws.Range("J6").Value = 0.2429203181515272
$ws.Range("M6").Value = 47.32925566666668
$ws.Range("N6").Value = 141.987767
$ws.Range("O6").Value = 0.3408416299313156
$ws.Range("P6").Value = 0.3408416299313156
$ws.Range("Q6").Value = 4726.485585283033
$ws.Range("R6").Value = 42538.3702675473
$ws.Range("S6").Value = 0.08279735718220027
$ws.Range("T6").Value = 0.0827973571822003
# Row 7
$ws.Range("I7").Value = 0.2429203181515272
$ws.Range("J7").Value = 0.2429203181515272
$ws.Range("M7").Value = 43.717953
$ws.Range("O7").Value = 0.3148348341399153
$ws.Range("P7").Value = 0.3148348341399154
$ws.Range("Q7").Value = 4365.846700143847
$ws.Range("R7").Value = 39292.62030129463
$ws.Range("S7").Value = 0.07647977807445151
$ws.Range("T7").Value = 0.07647977807445155
# Row 8
$ws.Range("I8").Value = 0.2429203181515272
$ws.Range("J8").Value = 0.2429203181515272
$ws.Range("M8").Value = 21.09134933333333
$ws.Range("N8").Value = 63.274048
$ws.Range("O8").Value = 0.1518893501062827
$ws.Range("P8").Value = 0.1518893501062827
$ws.Range("Q8").Value = 2106.265082642695
$ws.Range("R8").Value = 18956.38574378426
$ws.Range("S8").Value = 0.0368970092516469
$ws.Range("T8").Value = 0.03689700925164691
# Row 9
$ws.Range("I9").Value = 0.2429203181515272
$ws.Range("J9").Value = 0.2429203181515272
$ws.Range("M9").Value = 26.72140366666666
$ws.Range("N9").Value = 80.16421099999999
$ws.Range("O9").Value = 0.1924341858224863
$ws.Range("P9").Value = 0.1924341858224864
$ws.Range("Q9").Value = 2668.504447619684
$ws.Range("R9").Value = 24016.54002857716
$ws.Range("S9").Value = 0.04674617364322848
$ws.Range("T9").Value = 0.0467461736432285
# Row 10
$ws.Range("G10").Value = 16.49037766666667
$ws.Range("H10").Value = 49.471133
$ws.Range("I10").Value = 0.04011305879538658
$ws.Range("J10").Value = 0.04011305879538658
$ws.Range("M10").Value = 47.32925566666668
$ws.Range("N10").Value = 141.987767
$ws.Range("O10").Value = 0.3408416299313156
$ws.Range("P10").Value = 0.3408416299313156
$ws.Range("Q10").Value = 780.4773006255569
$ws.Range("R10").Value = 7024.295705630012
$ws.Range("S10").Value = 0.01367220034135025
$ws.Range("T10").Value = 0.01367220034135026
# Row 11
$ws.Range("G11").Value = 16.49037766666667
$ws.Range("H11").Value = 49.471133
$ws.Range("I11").Value = 0.04011305879538658
$ws.Range("J11").Value = 0.04011305879538658
$ws.Range("M11").Value = 43.717953
$ws.Range("O11").Value = 0.3148348341399153
$ws.Range("P11").Value = 0.3148348341399154
$ws.Range("Q11").Value = 720.925555783583
$ws.Range("R11").Value = 6488.330002052247
$ws.Range("S11").Value = 0.0126289882126902
$ws.Range("T11").Value = 0.01262898821269021
# Row 12
$ws.Range("G12").Value = 16.49037766666667
$ws.Range("H12").Value = 49.471133
$ws.Range("I12").Value = 0.04011305879538658
$ws.Range("J12").Value = 0.04011305879538658
$ws.Range("M12").Value = 21.09134933333333
$ws.Range("N12").Value = 63.274048
$ws.Range("O12").Value = 0.1518893501062827
$ws.Range("P12").Value = 0.1518893501062827
$ws.Range("Q12").Value = 347.8043160062649
$ws.Range("R12").Value = 3130.238844056384
$ws.Range("S12").Value = 0.006092746431206375
$ws.Range("T12").Value = 0.006092746431206376
# Row 13
$ws.Range("G13").Value = 16.49037766666667
$ws.Range("H13").Value = 49.471133
$ws.Range("I13").Value = 0.04011305879538658
$ws.Range("J13").Value = 0.04011305879538658
$ws.Range("M13").Value = 26.72140366666666
$ws.Range("N13").Value = 80.16421099999999
$ws.Range("O13").Value = 0.1924341858224863
$ws.Range("P13").Value = 0.1924341858224864
$ws.Range("Q13").Value = 440.6460382467848
$ws.Range("R13").Value = 3965.814344221063
$ws.Range("S13").Value = 0.00771912381013974
$ws.Range("T13").Value = 0.007719123810139741
# Row 14
$ws.Range("G14").Value = 249.4253923333333
$ws.Range("H14").Value = 748.276177
$ws.Range("I14").Value = 0.6067305206692575
$ws.Range("J14").Value = 0.6067305206692575
$ws.Range("M14").Value = 47.32925566666668
$ws.Range("N14").Value = 141.987767
$ws.Range("O14").Value = 0.3408416299313156
$ws.Range("P14").Value = 0.3408416299313156
$ws.Range("Q14").Value = 11805.11816350298
$ws.Range("R14").Value = 106246.0634715268
$ws.Range("S14").Value = 0.2067990195939855
$ws.Range("T14").Value = 0.2067990195939855
# Row 15
$ws.Range("G15").Value = 249.4253923333333
$ws.Range("H15").Value = 748.276177
$ws.Range("I15").Value = 0.6067305206692575
$ws.Range("J15").Value = 0.6067305206692575
$ws.Range("M15").Value = 43.717953
$ws.Range("O15").Value = 0.3148348341399153
$ws.Range("P15").Value = 0.3148348341399154
$ws.Range("Q15").Value = 10904.36757903523
$ws.Range("R15").Value = 98139.30821131705
$ws.Range("S15").Value = 0.1910199028425301
$ws.Range("T15").Value = 0.1910199028425302
# Row 16
$ws.Range("G16").Value = 249.4253923333333
$ws.Range("H16").Value = 748.276177
$ws.Range("I16").Value = 0.6067305206692575
$ws.Range("J16").Value = 0.6067305206692575
$ws.Range("M16").Value = 21.09134933333333
$ws.Range("N16").Value = 63.274048
$ws.Range("O16").Value = 0.1518893501062827
$ws.Range("P16").Value = 0.1518893501062827
$ws.Range("Q16").Value = 5260.718082306055
$ws.Range("R16").Value = 47346.46274075449
$ws.Range("S16").Value = 0.09215590447410005
$ws.Range("T16").Value = 0.09215590447410006
# Row 17
$ws.Range("G17").Value = 249.4253923333333
$ws.Range("H17").Value = 748.276177
$ws.Range("I17").Value = 0.6067305206692575
$ws.Range("J17").Value = 0.6067305206692575
$ws.Range("M17").Value = 26.72140366666666
$ws.Range("N17").Value = 80.16421099999999
$ws.Range("O17").Value = 0.1924341858224863
$ws.Range("P17").Value = 0.1924341858224864
$ws.Range("Q17").Value = 6664.996593255704
$ws.Range("R17").Value = 59984.96933930134
$ws.Range("S17").Value = 0.1167556937586418
$ws.Range("T17").Value = 0.1167556937586418
